$d = $word.ActiveDocument

$d.Content.Find.Execute("441×2=", $true, $false, $false, $false, $false, $true, 1, $false, "426×4=", 2) | Out-Null
$d.Content.Find.Execute("458×5=", $true, $false, $false, $false, $false, $true, 1, $false, "844×4=", 2) | Out-Null
$d.Content.Find.Execute("225×3=", $true, $false, $false, $false, $false, $true, 1, $false, "521×8=", 2) | Out-Null
$d.Content.Find.Execute("203×4=", $true, $false, $false, $false, $false, $true, 1, $false, "225×2=", 2) | Out-Null
$d.Content.Find.Execute("344×9=", $true, $false, $false, $false, $false, $true, 1, $false, "938×3=", 2) | Out-Null
$d.Content.Find.Execute("146×7=", $true, $false, $false, $false, $false, $true, 1, $false, "896×6=", 2) | Out-Null
$d.Content.Find.Execute("233×8=", $true, $false, $false, $false, $false, $true, 1, $false, "534×4=", 2) | Out-Null
$d.Content.Find.Execute("863×5=", $true, $false, $false, $false, $false, $true, 1, $false, "384×5=", 2) | Out-Null
$d.Content.Find.Execute("334×5=", $true, $false, $false, $false, $false, $true, 1, $false, "239×3=", 2) | Out-Null
$d.Content.Find.Execute("301×9=", $true, $false, $false, $false, $false, $true, 1, $false, "605×3=", 2) | Out-Null
$d.Content.Find.Execute("151×2=", $true, $false, $false, $false, $false, $true, 1, $false, "370×5=", 2) | Out-Null
$d.Content.Find.Execute("810×3=", $true, $false, $false, $false, $false, $true, 1, $false, "376×7=", 2) | Out-Null
$d.Content.Find.Execute("652×7=", $true, $false, $false, $false, $false, $true, 1, $false, "805×6=", 2) | Out-Null
$d.Content.Find.Execute("836×4=", $true, $false, $false, $false, $false, $true, 1, $false, "726×7=", 2) | Out-Null
$d.Content.Find.Execute("899×5=", $true, $false, $false, $false, $false, $true, 1, $false, "683×9=", 2) | Out-Null
$d.Content.Find.Execute("110×6=", $true, $false, $false, $false, $false, $true, 1, $false, "787×4=", 2) | Out-Null
$d.Content.Find.Execute("559×6=", $true, $false, $false, $false, $false, $true, 1, $false, "136×3=", 2) | Out-Null
$d.Content.Find.Execute("141×8=", $true, $false, $false, $false, $false, $true, 1, $false, "889×9=", 2) | Out-Null
$d.Content.Find.Execute("119×8=", $true, $false, $false, $false, $false, $true, 1, $false, "301×9=", 2) | Out-Null
$d.Content.Find.Execute("940×2=", $true, $false, $false, $false, $false, $true, 1, $false, "142×8=", 2) | Out-Null
$d.Content.Find.Execute("787×9=", $true, $false, $false, $false, $false, $true, 1, $false, "230×7=", 2) | Out-Null
$d.Content.Find.Execute("297×7=", $true, $false, $false, $false, $false, $true, 1, $false, "150×4=", 2) | Out-Null
$d.Content.Find.Execute("410×3=", $true, $false, $false, $false, $false, $true, 1, $false, "113×2=", 2) | Out-Null
$d.Content.Find.Execute("934×5=", $true, $false, $false, $false, $false, $true, 1, $false, "452×3=", 2) | Out-Null
$d.Content.Find.Execute("721×3=", $true, $false, $false, $false, $false, $true, 1, $false, "628×8=", 2) | Out-Null
